$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = '51.012.79'
$ws.Range("E2").Value = '  -0.42%  '

$ws.Range("D3").Value = '2.949.65'
$ws.Range("E3").Value = '  -0.34%  '

Set-TextValue "D4" '0.999'
$ws.Range("E4").Value = '  -0.02%  '

Set-TextValue "D5" '379.33'
$ws.Range("E5").Value = '  -0.55%  '

Set-TextValue "D6" '101.38'
$ws.Range("E6").Value = '  -1.13%  '

Set-TextValue "D7" '0.541'
$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  -1.36%  '

Set-TextValue "D10" '36.20'
$ws.Range("E10").Value = '  -0.99%  '

$ws.Range("E11").Value = '  -0.55%  '

Set-TextValue "D12" '0.0847'
$ws.Range("E12").Value = '  +0.92%  '

$ws.Range("D13").Value = '3.410.23'
$ws.Range("E13").Value = '  -0.46%  '

Set-TextValue "D14" '18.35'
$ws.Range("E14").Value = '  +1.66%  '

Set-TextValue "D15" '7.75'
$ws.Range("E15").Value = '  +4.50%  '

Set-TextValue "D16" '11.83'
$ws.Range("E16").Value = '  +65.77%  '

$ws.Range("D17").Value = '2.940.84'
$ws.Range("E17").Value = '  -0.38%  '

$ws.Range("E18").Value = '  +1.45%  '

$ws.Range("D19").Value = '50.965.06'
$ws.Range("E19").Value = '  -0.40%  '

Set-TextValue "D20" '3.06'

$ws.Range("E21").Value = '  -0.86%  '

$ws.Range("E22").Value = '  -0.46%  '

Set-TextValue "D23" '69.46'
$ws.Range("E23").Value = '  +1.37%  '

Set-TextValue "D24" '266.73'
$ws.Range("E24").Value = '  +1.73%  '

$ws.Range("E25").Value = '  +11.43%  '

Set-TextValue "D26" '8.14'
$ws.Range("E26").Value = '  -2.94%  '

$ws.Range("E27").Value = '  -0.07%  '

Set-TextValue "D28" '7.08'
$ws.Range("E28").Value = '  -7.91%  '

Set-TextValue "D29" '25.61'
$ws.Range("E29").Value = '  -0.35%  '

Set-TextValue "D30" '0.162'
$ws.Range("E30").Value = '  -4.22%  '

$ws.Range("E31").Value = '  -2.84%  '

$ws.Range("E32").Value = '  +2.89%  '

$ws.Range("B33").Value = 'Toncoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D33" '2.06'
$ws.Range("E33").Value = '  +0.23%  '

$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D34" '50.50'
$ws.Range("E34").Value = '  +0.16%  '

$ws.Range("E35").Value = '  -1.39%  '

Set-TextValue "D36" '0.0432'
$ws.Range("E36").Value = '  -5.70%  '

$ws.Range("E37").Value = '  -0.03%  '

Set-TextValue "D38" '3.09'
$ws.Range("E38").Value = '  +3.54%  '

$ws.Range("E39").Value = '  +0.60%  '

Set-TextValue "D40" '16.64'
$ws.Range("E40").Value = '  -1.00%  '

Set-TextValue "D41" '2.52'
$ws.Range("E41").Value = '  -0.72%  '

$ws.Range("E42").Value = '  +1.38%  '

Set-TextValue "D43" '117.69'
$ws.Range("E43").Value = '  -3.20%  '

$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D44" '3.52'
$ws.Range("E44").Value = '  +8.83%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D45" '21.34'
$ws.Range("E45").Value = '  +0.17%  '

Set-TextValue "D46" '2.02'
$ws.Range("E46").Value = '  -1.46%  '

$ws.Range("E47").Value = '  -2.07%  '

$ws.Range("D48").Value = '2.010.10'
$ws.Range("E48").Value = '  -0.05%  '

Set-TextValue "D49" '0.260'
$ws.Range("E49").Value = '  -4.73%  '

Set-TextValue "D50" '0.0313'
$ws.Range("E50").Value = '  -9.17%  '

Set-TextValue "D51" '5.28'
$ws.Range("E51").Value = '  +4.10%  '
